$wb = $excel.ActiveWorkbook

$wsInforme = $wb.Worksheets.Item("INFORME")
$wsPend    = $wb.Worksheets.Item("PENDIENTES")

# ---------------------------------------------------------------
# INFORME sheet (row 24): update two existing notes and add two
# new ones in E24 / F24, with partial rich-text coloring.
# ---------------------------------------------------------------

# C24 - plain text edit (DB -> BASE DE DATOS)
$wsInforme.Range("C24").Value = "Login Acceso (error corregido), (Modificar relaciones en Diagrama de clases(DE LA BASE DE DATOS) "

# D24 - append "(pendiente)" with "pendiente" colored dark red, trailing ")" default colored
$wsInforme.Range("D24").Value = "Cambio de diseño(formulario Usuario) Modelo para todos los formulario, Cambiar Filtro para buscar registro e implementar un nuevo componente(DataAnnotations)para validar datos.(pendiente)"
$rngD = $wsInforme.Range("D24")
$rngD.Characters(178, 9).Font.Color = 192
$rngD.Characters(187, 1).Font.Color = 0

# E24 - new cell: base text + "concluido" colored blue, trailing space default colored
$wsInforme.Range("E24").Value = "Formulario Usuario en C#(prueba de validacion y filtro al registro)concluido "
$rngE = $wsInforme.Range("E24")
$rngE.Characters(68, 9).Font.Color = 12611584
$rngE.Characters(77, 1).Font.Color = 0

# F24 - new cell: plain text, no rich text
$wsInforme.Range("F24").Value = "Formulario Empresa y vista empresa en  C#(validacion y filtro), se modifico el PROCEDIENTO ALMACENADO EN LA DB."

# ---------------------------------------------------------------
# View state: PENDIENTES becomes the active/selected tab, with a
# new selection; INFORME keeps its own (now different) selection.
# ---------------------------------------------------------------

$wsInforme.Range("F24").Select()

$wsPend.Activate()
$wsPend.Range("F26").Select()
